# Generate Report for Handback
#
# Refreshes the localization-status report after a handback: the
# "Ready for handoff" status becomes "Handed back: in sync with en-US",
# the Latest Handback DateTime stamps are refreshed, and the stale
# "handback file is not the latest" Error Detail messages are cleared
# now that the handback is in sync.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$handedBack = "Handed back: in sync with en-US"

# Overview sheet: zh-cn / de-de status columns
$ws1.Range("E2").Value = $handedBack
$ws1.Range("F2").Value = $handedBack

# zh-cn detail sheet
$ws2.Range("C2").Value = $handedBack
$ws2.Range("K2").Value = "2016-08-20 16:57:49"
$ws2.Range("P2").Value = ""

# de-de detail sheet
$ws3.Range("C2").Value = $handedBack
$ws3.Range("K2").Value = "2016-08-20 16:57:55"
$ws3.Range("P2").Value = ""

# Resize the Status / Error Detail columns to fit the new (longer / now
# blank) content, matching the widths the report generator lays out.
$ws1.Columns.Item(5).ColumnWidth = 29.09
$ws1.Columns.Item(6).ColumnWidth = 29.09
$ws2.Columns.Item(3).ColumnWidth = 29.09
$ws2.Columns.Item(16).ColumnWidth = 12.75
$ws3.Columns.Item(3).ColumnWidth = 29.09
$ws3.Columns.Item(16).ColumnWidth = 12.75
